$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---
# Copy the formatting (bold/border/centered style) from H1 so the new
# header cells match the look of the existing headers, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows 2-37: I = 1 (constant), J = same value as column H ---
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}

# --- Row 38 uses special values (not following the I=1 / J=H pattern) ---
$ws.Cells.Item(38, 9).Value2 = 3
$ws.Cells.Item(38, 10).Value2 = 4
